# Update the sign-off sheet table: the row that used to reference
# "Download the Bluetooth Spec Version 5.0" now references the new
# exercise "Program a simple Application".

$d = $word.ActiveDocument

$oldText = "Download the Bluetooth Spec Version 5.0"
$newText = "Program a simple Application"

$target = $null
$targetTable = $null

foreach ($tbl in $d.Tables) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $cell = $null
            $cell = $tbl.Cell($r, $c)
            if ($cell -ne $null -and $cell.Range.Text -like "*$oldText*") {
                $target = $cell
                $targetTable = $tbl
                break
            }
        }
        if ($target -ne $null) { break }
    }
    if ($target -ne $null) { break }
}

if ($target -ne $null) {
    # Trim the trailing cell-mark (and paragraph mark, if any) from the
    # cell's range so we only replace the visible text, leaving all of
    # the run/paragraph formatting (and rsid attributes) untouched.
    $rng = $target.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}
